$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.347.49"
$ws.Range("E2").Value = "  -1.65%  "

# Row 3
$ws.Range("D3").Value = "3.341.11"

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.22%  "

# Row 8
$ws.Range("D8").Value = "3.331.03"
$ws.Range("E8").Value = "  -3.50%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.626"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.15%  "

# Row 15
$ws.Range("D15").Value = "3.873.87"
$ws.Range("E15").Value = "  -3.52%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.07%  "

# Row 18
$ws.Range("D18").Value = "3.323.56"
$ws.Range("E18").Value = "  -4.09%  "

# Row 19
$ws.Range("D19").Value = "64.233.41"
$ws.Range("E19").Value = "  -1.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.70"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.977"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "442.68"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +8.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +12.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.56"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("E27").Value = "  -2.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.91%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.21%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "576.84"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.39%  "

# Row 34
$ws.Range("E34").Value = "  -2.63%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.36"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.56%  "

# Row 36
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("E37").Value = "  -8.56%  "

# Row 38
$ws.Range("E38").Value = "  -1.90%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.74%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0744"
$ws.Range("E40").Value = "  -5.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.33%  "

# Row 42
$ws.Range("D42").Value = "3.090.16"
$ws.Range("E42").Value = "  -3.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.14%  "

# Row 47
$ws.Range("E47").Value = "  -3.46%  "

# Row 48
$ws.Range("E48").Value = "  -1.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.60%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.66%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.83"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.56%  "
